$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.728.72'
$ws.Range('E2').Value = '  +3.73%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.010.69'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '565.57'
$ws.Range('E5').Value = '  +2.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.62'
$ws.Range('E6').Value = '  +7.45%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.001.16'
$ws.Range('E9').Value = '  +2.47%  '
$ws.Range('E10').Value = '  +5.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.37'
$ws.Range('E11').Value = '  +12.95%  '
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('E13').Value = '  +4.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.97'
$ws.Range('E14').Value = '  +3.19%  '
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.509.42'
$ws.Range('E16').Value = '  +2.39%  '
$ws.Range('E17').Value = '  +4.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.014.04'
$ws.Range('E18').Value = '  +2.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '59.738.66'
$ws.Range('E19').Value = '  +3.60%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '437.19'
$ws.Range('E20').Value = '  +4.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.66'
$ws.Range('E21').Value = '  +2.96%  '
$ws.Range('E22').Value = '  +4.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.10'
$ws.Range('E23').Value = '  +1.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.40'
$ws.Range('E24').Value = '  +2.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.56'
$ws.Range('E25').Value = '  +0.90%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.22'
$ws.Range('E27').Value = '  +11.03%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +3.63%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.84'
$ws.Range('E30').Value = '  +5.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.35'
$ws.Range('E31').Value = '  +5.94%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.97'
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.106'
$ws.Range('E33').Value = '  +9.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0₃0789'
$ws.Range('E34').Value = '  +13.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.996'
$ws.Range('E35').Value = '  +6.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.93'
$ws.Range('E36').Value = '  +4.42%  '
$ws.Range('E37').Value = '  +1.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.23'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.58'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.81'
$ws.Range('E40').Value = '  +9.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '402.82'
$ws.Range('E41').Value = '  +7.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0353'
$ws.Range('E42').Value = '  +2.24%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.765.16'
$ws.Range('E43').Value = '  +2.39%  '
$ws.Range('E44').Value = '  -0.89%  '
$ws.Range('E45').Value = '  +6.21%  '
$ws.Range('E46').Value = '  +0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '123.61'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.04'
$ws.Range('E48').Value = '  +3.94%  '
$ws.Range('E49').Value = '  +1.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.01'
$ws.Range('E50').Value = '  +19.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.65'
$ws.Range('E51').Value = '  +2.44%  '
